{"js": "// Office.js (Word JavaScript API) script\n// Applies the edits described by the diff:\n//  1. In the \"Conclusions\" paragraph, rewrite the dictionary-attack and\n//     brute-force sentences with concrete figures instead of the\n//     yellow-highlighted \"\u2206\" placeholders (some placeholders remain).\n//  2. Add two new paragraphs after \"Due to the size of my table there\u2026.\"\n//     describing a possible look-up-table optimisation.\n//  3. Re-type the \"What I learnt\" heading (same visible text).\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1. Locate the \"Conclusions\" paragraph that contains the sentence we\n//    need to rewrite, by searching for a stable, unique anchor phrase.\n// ---------------------------------------------------------------------\nconst targetSearch = body.search(\"Using a dictionary attack naturally cracked\", { matchCase: true });\ntargetSearch.load(\"items\");\nawait context.sync();\n\nif (targetSearch.items.length === 0) {\n  throw new Error(\"Could not find the target 'Conclusions' paragraph.\");\n}\n\nconst conclusionsPara = targetSearch.items[0].paragraphs.getFirst();\nconclusionsPara.load(\"text\");\nawait context.sync();\n\n// Rebuild the whole paragraph via OOXML so the final run/highlight layout\n// matches the target precisely: the three remaining \"\u2206\" placeholders stay\n// yellow-highlighted, while the newly supplied figures (\"969 of the 100\n// thousand\", \"3:43\", \"3:43\" again and \"18\") are plain (unhighlighted) text.\nconst conclusionsOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData>` +\n  `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body>` +\n  `<w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr>` +\n  `<w:r><w:t xml:space=\"preserve\">The generation took </w:t></w:r>` +\n  `<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>\\u2206</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\"> and created a \\u2206Mb file of 500,000 rows but theoretically the table holds 5,000,000,000 hashes. Of the 100,000 top passwords, it was able to crack </w:t></w:r>` +\n  `<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>\\u2206</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">, the longest took </w:t></w:r>` +\n  `<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>\\u2206</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\"> and the shortest too </w:t></w:r>` +\n  `<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>\\u2206</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">. Using a dictionary attack cracked </w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">969 of the 100 thousand </w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">in a time of </w:t></w:r>` +\n  `<w:r><w:t>3:43</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">, there is never a guarantee that </w:t></w:r>` +\n  `<w:r><w:t>the sought password would be known to check against</w:t></w:r>` +\n  `<w:r><w:t>. Running the brute</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">force attack for </w:t></w:r>` +\n  `<w:r><w:t>3:4</w:t></w:r>` +\n  `<w:r><w:t>3</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\"> cracked </w:t></w:r>` +\n  `<w:r><w:t>18</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\"> of the passwords</w:t></w:r>` +\n  `<w:r><w:t>.</w:t></w:r>` +\n  `</w:p>` +\n  `</w:body></w:document>` +\n  `</pkg:xmlData></pkg:part></pkg:package>`;\n\nconclusionsPara.insertOoxml(conclusionsOoxml, \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2. Insert the two new paragraphs right after\n//    \"Due to the size of my table there\u2026.\"\n// ---------------------------------------------------------------------\nconst dueSearch = body.search(\"Due to the size of my table there\", { matchCase: true });\ndueSearch.load(\"items\");\nawait context.sync();\n\nif (dueSearch.items.length === 0) {\n  throw new Error(\"Could not find the 'Due to the size of my table' paragraph.\");\n}\n\nconst duePara = dueSearch.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// First, a blank paragraph (same justification as its neighbours).\nconst blankPara = duePara.insertParagraph(\"\", \"After\");\nblankPara.alignment = Word.Alignment.justified;\n\n// Then the paragraph with the actual note text.\nconst notePara = blankPara.insertParagraph(\n  \"Currently _crack generates a hash then runs through all uncracked hashes to compare. Would use a look up table if I was to do it again, performace drasticly dropped for large uncriokded sets\",\n  \"After\"\n);\nnotePara.alignment = Word.Alignment.justified;\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3. Re-type the \"What I learnt\" heading (text is unchanged, but the\n//    original commit re-entered it, visible in the OOXML as a run split).\n// ---------------------------------------------------------------------\nconst headingSearch = body.search(\"What I learnt\", { matchCase: true });\nheadingSearch.load(\"items\");\nawait context.sync();\n\nif (headingSearch.items.length > 0) {\n  headingSearch.items[0].insertText(\"What I learnt\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the edits described by the diff:\n#  1. In the \"Conclusions\" paragraph, rewrite the dictionary-attack and\n#     brute-force sentences with concrete figures instead of the\n#     yellow-highlighted \"\u2206\" placeholders (some placeholders remain).\n#  2. Add two new paragraphs after \"Due to the size of my table there\u2026.\"\n#     describing a possible look-up-table optimisation.\n#  3. Re-type the \"What I learnt\" heading (same visible text).\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($doc, $findText, $replaceText, $clearHighlight) {\n    $rng = $doc.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    if ($clearHighlight) {\n        # wdNoHighlight = 0 -- remove the yellow highlight from the\n        # text that gets substituted in.\n        $find.Replacement.Highlight = 0\n    }\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 1. Rewrite the two sentences in the \"Conclusions\" paragraph.\n# ---------------------------------------------------------------------\nReplace-DocText $d \". Using a dictionary attack naturally cracked all the passwords in a time of \u2206, but there is never a guarantee that \" \". Using a dictionary attack cracked 969 of the 100 thousand in a time of 3:43, there is never a guarantee that \" $true\nReplace-DocText $d \"force attack for \u2206 cracked \u2206 of the passwords. \" \"force attack for 3:43 cracked 18 of the passwords.\" $true\n\n# ---------------------------------------------------------------------\n# 2. Insert the two new paragraphs right after\n#    \"Due to the size of my table there\u2026.\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Due to the size of my table there\"\n$find.Execute() | Out-Null\n$rng.Expand(4) | Out-Null   # wdParagraph -- grow to the whole paragraph\n\n# Blank paragraph (inherits the \"justify\" alignment of its neighbour).\n$rng.InsertParagraphAfter()\n$rng.Collapse(0) | Out-Null\n\n# Paragraph with the actual note text.\n$rng.InsertParagraphAfter()\n$rng.Collapse(0) | Out-Null\n$rng.Move(1, 1) | Out-Null\n$rng.Text = \"Currently _crack generates a hash then runs through all uncracked hashes to compare. Would use a look up table if I was to do it again, performace drasticly dropped for large uncriokded sets\"\n\n# ---------------------------------------------------------------------\n# 3. Re-type the \"What I learnt\" heading (text is unchanged, but the\n#    original commit re-entered it, visible in the OOXML as a run split).\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"What I learnt\"\nif ($find.Execute()) {\n    $rng.Text = \"What I learnt\"\n}\n"}
